$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Remove the GREECE - SUPER LEAGUE row (row 6); remaining rows shift up.
$ws.Rows(6).Delete()

# Step 2: Remove the UKRAINE - PREMIER LEAGUE row, which is now at row 10
# (it was row 11 before the previous delete).
$ws.Rows(10).Delete()

# Step 3: Insert a new blank row at position 8 for the new ROMANIA - LIGA 1 match,
# pushing the Turkey rows back down to 9 and 10.
$ws.Rows(8).Insert()

# Step 4: Populate the new row 8 with the ROMANIA - LIGA 1 match data.
$ws.Cells.Item(8, 1).Value = "feF9vkFc"
$ws.Cells.Item(8, 2).Value = "25/11/2024"
$ws.Cells.Item(8, 3).Value = "15:30"
$ws.Cells.Item(8, 4).Value = "ROMANIA - LIGA 1"
$ws.Cells.Item(8, 5).Value = "UTA Arad"
$ws.Cells.Item(8, 6).Value = "Univ. Craiova"
$ws.Cells.Item(8, 7).Value = 4
$ws.Cells.Item(8, 8).Value = 3.1
$ws.Cells.Item(8, 9).Value = 1.95
$ws.Cells.Item(8, 10).Value = 4.33
$ws.Cells.Item(8, 11).Value = 2.1
$ws.Cells.Item(8, 12).Value = 2.63
$ws.Cells.Item(8, 13).Value = 1.06
$ws.Cells.Item(8, 14).Value = 10
$ws.Cells.Item(8, 15).Value = 1.33
$ws.Cells.Item(8, 16).Value = 3.25
$ws.Cells.Item(8, 17).Value = 2.05
$ws.Cells.Item(8, 18).Value = 1.72
$ws.Cells.Item(8, 19).Value = 1.44
$ws.Cells.Item(8, 20).Value = 2.63
$ws.Cells.Item(8, 21).Value = 1.8
$ws.Cells.Item(8, 22).Value = 1.8
$ws.Cells.Item(8, 23).Value = 11
$ws.Cells.Item(8, 24).Value = 21
$ws.Cells.Item(8, 25).Value = 15
$ws.Cells.Item(8, 26).Value = 41
$ws.Cells.Item(8, 27).Value = 34
$ws.Cells.Item(8, 28).Value = 41
$ws.Cells.Item(8, 29).Value = 8.5
$ws.Cells.Item(8, 30).Value = 6
$ws.Cells.Item(8, 31).Value = 15
$ws.Cells.Item(8, 32).Value = 51
$ws.Cells.Item(8, 33).Value = 251
$ws.Cells.Item(8, 34).Value = 7
$ws.Cells.Item(8, 35).Value = 9
$ws.Cells.Item(8, 36).Value = 9
$ws.Cells.Item(8, 37).Value = 17
$ws.Cells.Item(8, 38).Value = 17
$ws.Cells.Item(8, 39).Value = 29
$ws.Cells.Item(8, 40).Value = 5.5
$ws.Cells.Item(8, 41).Value = 21
$ws.Cells.Item(8, 42).Value = 29
$ws.Cells.Item(8, 43).Value = 67
$ws.Cells.Item(8, 44).Value = 101
$ws.Cells.Item(8, 45).Value = 251
$ws.Cells.Item(8, 46).Value = 2.63
$ws.Cells.Item(8, 47).Value = 8.5
$ws.Cells.Item(8, 48).Value = 51
$ws.Cells.Item(8, 49).Value = 4
$ws.Cells.Item(8, 50).Value = 11
$ws.Cells.Item(8, 51).Value = 23
$ws.Cells.Item(8, 52).Value = 41
$ws.Cells.Item(8, 53).Value = 51
$ws.Cells.Item(8, 54).Value = 151
$ws.Cells.Item(8, 55).Value = 51
$ws.Cells.Item(8, 56).Value = 51

# Step 5: Apply updated odds values to the other rows.
# Row 2 tweaks (from before row 2)
$ws.Cells.Item(2, 11).Value = 1.91
# Row 3 tweaks (from before row 3)
$ws.Cells.Item(3, 7).Value = 1.24
$ws.Cells.Item(3, 8).Value = 5.5
$ws.Cells.Item(3, 9).Value = 12
$ws.Cells.Item(3, 21).Value = 2.02
$ws.Cells.Item(3, 22).Value = 1.72
$ws.Cells.Item(3, 23).Value = 7
$ws.Cells.Item(3, 24).Value = 6.7
$ws.Cells.Item(3, 25).Value = 9.5
$ws.Cells.Item(3, 26).Value = 7.9
$ws.Cells.Item(3, 28).Value = 30
$ws.Cells.Item(3, 34).Value = 25
$ws.Cells.Item(3, 35).Value = 110
$ws.Cells.Item(3, 42).Value = 15
$ws.Cells.Item(3, 49).Value = 11.25
$ws.Cells.Item(3, 50).Value = 70
$ws.Cells.Item(3, 53).Value = 500
# Row 4 tweaks (from before row 4)
$ws.Cells.Item(4, 13).Value = 1.02
$ws.Cells.Item(4, 15).Value = 1.13
# Row 5 tweaks (from before row 5)
$ws.Cells.Item(5, 8).Value = 4.33
$ws.Cells.Item(5, 11).Value = 2.5
$ws.Cells.Item(5, 12).Value = 1.95
$ws.Cells.Item(5, 15).Value = 1.2
$ws.Cells.Item(5, 16).Value = 4.33
$ws.Cells.Item(5, 17).Value = 1.67
$ws.Cells.Item(5, 18).Value = 2.15
$ws.Cells.Item(5, 21).Value = 1.83
$ws.Cells.Item(5, 22).Value = 1.83
$ws.Cells.Item(5, 28).Value = 41
# Row 6 tweaks (from before row 7)
$ws.Cells.Item(6, 13).Value = 1.05
$ws.Cells.Item(6, 15).Value = 1.29
# Row 7 tweaks (from before row 8)
$ws.Cells.Item(7, 21).Value = 1.63
# Row 9 tweaks (from before row 9)
$ws.Cells.Item(9, 7).Value = 1.5
$ws.Cells.Item(9, 8).Value = 4.5
$ws.Cells.Item(9, 9).Value = 5.5
$ws.Cells.Item(9, 10).Value = 2
$ws.Cells.Item(9, 11).Value = 2.5
$ws.Cells.Item(9, 13).Value = 1.03
$ws.Cells.Item(9, 14).Value = 17
$ws.Cells.Item(9, 17).Value = 1.53
$ws.Cells.Item(9, 18).Value = 2.4
$ws.Cells.Item(9, 24).Value = 9
$ws.Cells.Item(9, 26).Value = 12
$ws.Cells.Item(9, 29).Value = 17
$ws.Cells.Item(9, 34).Value = 19
$ws.Cells.Item(9, 36).Value = 17
$ws.Cells.Item(9, 37).Value = 51
$ws.Cells.Item(9, 41).Value = 7.5
$ws.Cells.Item(9, 43).Value = 21
$ws.Cells.Item(9, 50).Value = 26
$ws.Cells.Item(9, 53).Value = 81
# Row 10 tweaks (from before row 10)
$ws.Cells.Item(10, 12).Value = 4.75
$ws.Cells.Item(10, 15).Value = 1.44
$ws.Cells.Item(10, 16).Value = 2.63
$ws.Cells.Item(10, 19).Value = 1.53
$ws.Cells.Item(10, 20).Value = 2.38
$ws.Cells.Item(10, 46).Value = 2.38
Write-Host "Edit complete."
